$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.569.36'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '1.597.75'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = '208.30'
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("E6").Value = '  -3.46%  '
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("D8").Value = '22.33'
$ws.Range("E8").Value = '  -4.41%  '
$ws.Range("D9").Value = '0.252'
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").Value = '0.0592'
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").Value = '1.824.80'
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("D13").Value = '1.591.32'
$ws.Range("E13").Value = '  -2.22%  '
$ws.Range("D14").Value = '3.87'
$ws.Range("E14").Value = '  -3.79%  '
$ws.Range("D15").Value = '0.540'
$ws.Range("E15").Value = '  -3.76%  '
$ws.Range("D16").Value = '63.45'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").Value = '27.559.85'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").Value = '218.93'
$ws.Range("E18").Value = '  -4.82%  '
$ws.Range("E19").Value = '  -3.07%  '
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '4.21'
$ws.Range("E22").Value = '  -2.48%  '
$ws.Range("D23").Value = '9.70'
$ws.Range("E23").Value = '  -3.94%  '
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("D25").Value = '154.59'
$ws.Range("D26").Value = '6.75'
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").Value = '15.07'
$ws.Range("E28").Value = '  -2.92%  '
$ws.Range("D29").Value = '0.106'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("E32").Value = '  -4.45%  '
$ws.Range("D33").Value = '1.365.32'
$ws.Range("E33").Value = '  -2.34%  '
$ws.Range("D34").Value = '2.95'
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("E35").Value = '  -2.71%  '
$ws.Range("D36").Value = '0.978'
$ws.Range("E36").Value = '  -4.28%  '
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("E38").Value = '  -2.52%  '
$ws.Range("D39").Value = '0.539'
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").Value = '0.815'
$ws.Range("E40").Value = '  -4.35%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").Value = '0.981'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").Value = '5.36'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("E44").Value = '  -3.60%  '
$ws.Range("D45").Value = '64.15'
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").Value = '1.734.91'
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("D47").Value = '2.10'
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("D48").Value = '88.15'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = '0.0973'
$ws.Range("E49").Value = '  -4.09%  '
$ws.Range("D50").Value = '0.0₇0982'
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("D51").Value = '0.0498'
$ws.Range("E51").Value = '  -1.01%  '
